# Auto-generated edit script: updates cached market-price / profit values
# on the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR),
# as produced by the scheduled price-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 327.2
$ws.Range("I8").Value = 557.5
$ws.Range("K8").Value = 1672.5
$ws.Range("M8").Value = -1533.5
$ws.Range("H9").Value = 83
$ws.Range("I9").Value = 99.5
$ws.Range("K9").Value = 99.5
$ws.Range("M9").Value = 69.5
$ws.Range("H15").Value = 493.56097
$ws.Range("I15").Value = 493.56097
$ws.Range("K15").Value = 1480.68291
$ws.Range("M15").Value = -1311.68291
$ws.Range("H33").Value = 126.72727
$ws.Range("I33").Value = 115.5
$ws.Range("K33").Value = 115.5
$ws.Range("M33").Value = 113.5
$ws.Range("H39").Value = 538.44446
$ws.Range("I39").Value = 555.75
$ws.Range("K39").Value = 1667.25
$ws.Range("M39").Value = -1371.25
$ws.Range("H62").Value = 6520.1177
$ws.Range("J62").Value = 9385.200000000001
$ws.Range("L62").Value = 9385.200000000001
$ws.Range("N62").Value = -10633.2
$ws.Range("H65").Value = 6520.1177
$ws.Range("J65").Value = 9385.200000000001
$ws.Range("L65").Value = 46926
$ws.Range("N65").Value = -53166
$ws.Range("H70").Value = 15884249
$ws.Range("I70").Value = 27779186
$ws.Range("K70").Value = 83337558
$ws.Range("M70").Value = -83337288
$ws.Range("H73").Value = 15884249
$ws.Range("I73").Value = 27779186
$ws.Range("K73").Value = 83337558
$ws.Range("M73").Value = -83336622
$ws.Range("H100").Value = 5854.727
$ws.Range("I100").Value = 4933.6665
$ws.Range("J100").Value = 9999.5
$ws.Range("K100").Value = 4933.6665
$ws.Range("L100").Value = 9999.5
$ws.Range("M100").Value = -4392.6665
$ws.Range("N100").Value = -11081.5
$ws.Range("H107").Value = 836.0476
$ws.Range("I107").Value = 913.2105
$ws.Range("J107").Value = 103
$ws.Range("K107").Value = 913.2105
$ws.Range("L107").Value = 103
$ws.Range("M107").Value = 1006.7895
$ws.Range("N107").Value = -3943
$ws.Range("H111").Value = 1422.7222
$ws.Range("I111").Value = 1180.8334
$ws.Range("J111").Value = 1906.5
$ws.Range("K111").Value = 3542.5002
$ws.Range("L111").Value = 5719.5
$ws.Range("M111").Value = -475.5001999999999
$ws.Range("N111").Value = -11853.5
$ws.Range("H112").Value = 11722.333
$ws.Range("J112").Value = 13866.9
$ws.Range("L112").Value = 41600.7
$ws.Range("N112").Value = -43816.7
$ws.Range("H116").Value = 7980.8057
$ws.Range("I116").Value = 8363.666999999999
$ws.Range("K116").Value = 8363.666999999999
$ws.Range("M116").Value = -4921.666999999999
$ws.Range("H137").Value = 18555.684
$ws.Range("I137").Value = 1994.8572
$ws.Range("J137").Value = 28216.166
$ws.Range("K137").Value = 5984.571599999999
$ws.Range("L137").Value = 84648.49800000001
$ws.Range("M137").Value = -3434.571599999999
$ws.Range("N137").Value = -89748.49800000001
$ws.Range("H138").Value = 1394309.1
$ws.Range("I138").Value = 2243.842
$ws.Range("J138").Value = 2109153.5
$ws.Range("K138").Value = 6731.526
$ws.Range("L138").Value = 6327460.5
$ws.Range("M138").Value = -1591.526
$ws.Range("N138").Value = -6337740.5
$ws.Range("H141").Value = 1963.8334
$ws.Range("I141").Value = 1963.8334
$ws.Range("K141").Value = 5891.5002
$ws.Range("M141").Value = -711.5002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3982.158
$ws.Range("I63").Value = 2406.7778
$ws.Range("K63").Value = 2406.7778
$ws.Range("M63").Value = -1720.7778
$ws.Range("H66").Value = 3982.158
$ws.Range("I66").Value = 2406.7778
$ws.Range("K66").Value = 12033.889
$ws.Range("M66").Value = -8601.888999999999
$ws.Range("H74").Value = 3407.5095
$ws.Range("I74").Value = 1230.8422
$ws.Range("J74").Value = 8921.733
$ws.Range("K74").Value = 1230.8422
$ws.Range("L74").Value = 8921.733
$ws.Range("M74").Value = -356.8422
$ws.Range("N74").Value = -10669.733
$ws.Range("H77").Value = 3407.5095
$ws.Range("I77").Value = 1230.8422
$ws.Range("J77").Value = 8921.733
$ws.Range("K77").Value = 6154.211
$ws.Range("L77").Value = 44608.665
$ws.Range("M77").Value = -1786.211
$ws.Range("N77").Value = -53344.665
$ws.Range("H97").Value = 799.70966
$ws.Range("I97").Value = 693.05
$ws.Range("J97").Value = 993.63635
$ws.Range("K97").Value = 693.05
$ws.Range("L97").Value = 993.63635
$ws.Range("M97").Value = -197.05
$ws.Range("N97").Value = -1985.63635
$ws.Range("H102").Value = 670.9167
$ws.Range("I102").Value = 450.375
$ws.Range("K102").Value = 450.375
$ws.Range("M102").Value = 1171.625
$ws.Range("H110").Value = 23050.264
$ws.Range("I110").Value = 24711.715
$ws.Range("J110").Value = 3666.6667
$ws.Range("K110").Value = 24711.715
$ws.Range("L110").Value = 3666.6667
$ws.Range("M110").Value = -22666.715
$ws.Range("N110").Value = -7756.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3572041
$ws.Range("I94").Value = 693.15
$ws.Range("K94").Value = 693.15
$ws.Range("M94").Value = -242.15
$ws.Range("H105").Value = 2784.85
$ws.Range("I105").Value = 3057.375
$ws.Range("K105").Value = 3057.375
$ws.Range("M105").Value = -1310.375
$ws.Range("H107").Value = 882.0741
$ws.Range("I107").Value = 757.7895
$ws.Range("K107").Value = 757.7895
$ws.Range("M107").Value = 1162.2105

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1756545.8
$ws.Range("I31").Value = 3705791.5
$ws.Range("K31").Value = 3705791.5
$ws.Range("M31").Value = -3705496.5
$ws.Range("H34").Value = 1756545.8
$ws.Range("I34").Value = 3705791.5
$ws.Range("K34").Value = 3705791.5
$ws.Range("M34").Value = -3705589.5
$ws.Range("H99").Value = 7529.8667
$ws.Range("I99").Value = 7440.4443
$ws.Range("K99").Value = 7440.4443
$ws.Range("M99").Value = -5942.4443
$ws.Range("H107").Value = 550.5
$ws.Range("I107").Value = 357.2857
$ws.Range("K107").Value = 357.2857
$ws.Range("M107").Value = 1562.7143
$ws.Range("H126").Value = 7529.8667
$ws.Range("I126").Value = 7440.4443
$ws.Range("K126").Value = 22321.3329
$ws.Range("M126").Value = -19851.3329
$ws.Range("H134").Value = 2166.1292
$ws.Range("I134").Value = 1573.25
$ws.Range("J134").Value = 4198.857
$ws.Range("K134").Value = 4719.75
$ws.Range("L134").Value = 12596.571
$ws.Range("M134").Value = -2184.75
$ws.Range("N134").Value = -17666.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 157.65218
$ws.Range("J12").Value = 184
$ws.Range("L12").Value = 552
$ws.Range("N12").Value = -898
$ws.Range("H68").Value = 4031.074
$ws.Range("I68").Value = 1508
$ws.Range("J68").Value = 4751.952
$ws.Range("K68").Value = 4524
$ws.Range("L68").Value = 14255.856
$ws.Range("M68").Value = -3713
$ws.Range("N68").Value = -15877.856
$ws.Range("H71").Value = 4031.074
$ws.Range("I71").Value = 1508
$ws.Range("J71").Value = 4751.952
$ws.Range("K71").Value = 13572
$ws.Range("L71").Value = 42767.568
$ws.Range("M71").Value = -9516
$ws.Range("N71").Value = -50879.568
$ws.Range("H132").Value = 1548.125
$ws.Range("I132").Value = 1586.2
$ws.Range("J132").Value = 1530.8182
$ws.Range("K132").Value = 14275.8
$ws.Range("L132").Value = 13777.3638
$ws.Range("M132").Value = -11745.8
$ws.Range("N132").Value = -18837.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 34999.8
$ws.Range("J63").Value = 34999.8
$ws.Range("L63").Value = 34999.8
$ws.Range("N63").Value = -36371.8
$ws.Range("H66").Value = 34999.8
$ws.Range("J66").Value = 34999.8
$ws.Range("L66").Value = 104999.4
$ws.Range("N66").Value = -111863.4
$ws.Range("H113").Value = 1227.3334
$ws.Range("I113").Value = 943.25
$ws.Range("K113").Value = 943.25
$ws.Range("M113").Value = 1226.75
$ws.Range("H132").Value = 2214.074
$ws.Range("I132").Value = 1991.24
$ws.Range("K132").Value = 5973.72
$ws.Range("M132").Value = -3443.72
$ws.Range("N136").ClearContents()
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2121.111
$ws.Range("I22").Value = 2070
$ws.Range("J22").Value = 2135.7144
$ws.Range("K22").Value = 2070
$ws.Range("L22").Value = 2135.7144
$ws.Range("M22").Value = -1775
$ws.Range("N22").Value = -2725.7144
$ws.Range("H27").Value = 2121.111
$ws.Range("I27").Value = 2070
$ws.Range("J27").Value = 2135.7144
$ws.Range("K27").Value = 2070
$ws.Range("L27").Value = 2135.7144
$ws.Range("M27").Value = -1963
$ws.Range("N27").Value = -2349.7144
$ws.Range("H61").Value = 2447
$ws.Range("I61").Value = 2447
$ws.Range("K61").Value = 2447
$ws.Range("M61").Value = -2245
$ws.Range("H113").Value = 2447
$ws.Range("I113").Value = 2447
$ws.Range("K113").Value = 2447
$ws.Range("M113").Value = -277
$ws.Range("H122").Value = 4425
$ws.Range("I122").Value = 4000
$ws.Range("J122").Value = 4850
$ws.Range("K122").Value = 12000
$ws.Range("L122").Value = 14550
$ws.Range("M122").Value = -9550
$ws.Range("N122").Value = -19450
$ws.Range("H136").Value = 4991.091
$ws.Range("I136").Value = 2818.1667
$ws.Range("J136").Value = 7598.6
$ws.Range("K136").Value = 8454.500100000001
$ws.Range("L136").Value = 22795.8
$ws.Range("M136").Value = -5904.500100000001
$ws.Range("N136").Value = -27895.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 613.7447
$ws.Range("I113").Value = 659.32434
$ws.Range("J113").Value = 445.1
$ws.Range("K113").Value = 1977.97302
$ws.Range("L113").Value = 1335.3
$ws.Range("M113").Value = 192.0269800000001
$ws.Range("N113").Value = -5675.3
$ws.Range("N137").ClearContents()
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0

